$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns in the sheet that receive new values (B,C,D,F,G,H,J,L,M,O)
# Columns E,I,K,N remain 0 and are not touched.
$cols = @(2,3,4,6,7,8,10,12,13,15)

$data = @(
        @(16.20463809168008,12.39663578697463,9.933165352441552,31.44008403682198,31.49775227646713,15.09079098713568,10.88435774613836,11.47108939899972,16.07780712514151,23.31350193068468),
        @(15.65986545921441,12.26094000174127,9.944289296897605,31.57553455690758,31.65627782474158,15.15344406643017,10.91254921279299,11.46729675016428,15.94603786500474,23.42482096507096),
        @(15.3160916985288,12.17740242531245,9.952293000589563,31.6665004871389,31.76453091899514,15.19455817276572,10.93076212961629,11.46605790625599,15.86586190064791,23.49862279337472),
        @(15.17385940297495,12.14333265935346,9.955850172767185,31.70552641828007,31.81137646322028,15.21197774313494,10.93841184061089,11.46582840246679,15.833397839712,23.53006622898775),
        @(15.15011845008282,12.13767454993651,9.956458705806947,31.71212468405902,31.81931969701037,15.21491043349943,10.93969584760235,11.46580695923232,15.82802049623296,23.53536998214346),
        @(15.3141819087396,12.17694302329554,9.952339776351485,31.66701888838815,31.76515165306769,15.1947904054309,10.93086437304958,11.46605369462035,15.86542320350539,23.49904131211526),
        @(16.01883655855065,12.34991012219656,9.93675762390234,31.48516607388518,31.55013769050198,15.11184492659168,10.89389110122557,11.46955627988463,16.03223619251143,23.35075193710646),
        @(17.3194763281551,12.68614351645142,9.91549019356456,31.19062986125462,31.21572576800669,14.97016726665785,10.82852266787361,11.48501459573887,16.36406626646356,23.10331153104724),
        @(18.21666952688174,12.92977214879741,9.905496621868712,31.0123413529969,31.02400884879104,14.87885484323626,10.78480353915019,11.50152863457201,16.60926348695825,22.94808297781518),
        @(18.61059237752489,13.03953462231194,9.902165792417259,30.93956240468923,30.94867376623441,14.84008724372908,10.76584061342061,11.51014269620052,16.7208003187817,22.8832634788292),
        @(18.7576067784601,13.08091849990799,9.90107856010416,30.91320446505748,30.92186689653133,14.82580537679083,10.75879218837945,11.51356133964831,16.76300937428518,22.8595538502982),
        @(18.72604210227174,13.07201421365474,9.901304983984485,30.91882759241276,30.92756349998658,14.82886350302301,10.76030431288898,11.51281813355234,16.75392054822919,22.86462291173277),
        @(18.62273113386575,13.04294307757468,9.902072860353407,30.93736980755144,30.94643378814354,14.838904276941,10.76525808475443,11.5104208201006,16.72427358780161,22.88129610032214),
        @(18.55916637593266,13.02511181669313,9.902565856417954,30.94888410518843,30.95821686697722,14.84510645225499,10.76830964052858,11.50897274841457,16.70610959098666,22.89161787991876),
        @(18.19062985465057,12.92257510034909,9.905738697280983,31.01726556047019,31.02917230032866,14.88144416508918,10.78606137468321,11.50098770787746,16.60197189914927,22.95243593159064),
        @(17.96082012037771,12.85937942468274,9.907995929691021,31.06135107977923,31.07575271285483,14.90444595163324,10.79718800062911,11.49637010958713,16.53806644427051,22.99123202346934),
        @(17.82730743095421,12.82293227064783,9.909408662635496,31.08749125758406,31.10366180498315,14.9179367699904,10.80367485476418,11.49381799946457,16.50131087952338,23.01409176883912),
        @(17.78187714462838,12.81057579516395,9.909906663616615,31.09647626918548,31.1133028349617,14.92254931808549,10.80588617252319,11.49297178026459,16.48886709822451,23.02192524849903),
        @(17.98542248879658,12.86611712653719,9.907743804892274,31.05657699246563,31.07067842003098,14.90197037955609,10.79599454083875,11.49685092923329,16.54486936393007,22.98704566044553),
        @(18.65313539495887,13.05148709779437,9.901842597557804,30.93189085447451,30.940844323032,14.83594423970286,10.76379945276524,11.51112073104563,16.73298257859773,22.87637606633082),
        @(19.07691724822235,13.17157103777662,9.899000148651153,30.85740873715679,30.86602644393568,14.79511565492436,10.74352968800337,11.52135929199765,16.85575274013224,22.80892188740604),
        @(18.85192386688809,13.10758644787386,9.900424640159224,30.89651854850053,30.90503574632679,14.81669400662406,10.75427764232009,11.51581188800195,16.79025260657011,22.84447647164496),
        @(17.9743040852383,12.86307136799791,9.907857432197975,31.05873288186586,31.07296899081492,14.90308875467539,10.79653382359097,11.4966332309288,16.54179381007721,22.98893658547605),
        @(16.97729299327645,12.59566644224938,9.920252215327114,31.26363634370494,31.29677040309819,15.00625015400289,10.84544709157884,11.47992105801872,16.27394441195466,23.1655953368332)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
